$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.241.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "'1.656.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").Value = "'219.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'0.5232"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.2667"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.06365"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'20.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'0.07723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "'4.598"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "'1.639.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "'1.884.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "'0.5653"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'0.0₅8277"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'65.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "'26.227.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "'4.699"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'10.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'192.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("D23").Value = "'6.013"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'143.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").Value = "'0.1199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "'7.285"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'1.504"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'0.05637"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").Value = "'1.279"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").Value = "'3.505"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "'3.358"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "'2.805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9460"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").Value = "'0.5755"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "'5.917"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'0.8472"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'1.019.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.65%  "
$ws.Range("D45").Value = "'101.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").Value = "'1.795.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'58.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'0.05315"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("D51").Value = "'0.4350"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.40%  "
